$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update unit price values in column D to their new figures
$ws.Range("D17").Value = 325
$ws.Range("D20").Value = 278
$ws.Range("D22").Value = 164
$ws.Range("D29").Value = 127

# Reflect the updated view/selection state (scrolled up, D17 now selected)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D17").Select()
